$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 holds the most recent month (2025-08); update raw inputs
$ws.Range("B21").Value = 6261
$ws.Range("D21").Value = 5681666

# Recompute dependent metrics for row 21, matching the pattern used
# by the other rows in the sheet (row 21 vs. row 9, 12 months prior)
$b21 = $ws.Range("B21").Value()
$c21 = $ws.Range("C21").Value()
$d21 = $ws.Range("D21").Value()

$b9 = $ws.Range("B9").Value()
$c9 = $ws.Range("C9").Value()
$d9 = $ws.Range("D9").Value()

$ws.Range("E21").Value = $d21 / $b21
$ws.Range("F21").Value = ($b21 - $b9) / $b9 * 100
$ws.Range("G21").Value = ($c21 - $c9) / $c9 * 100
$ws.Range("H21").Value = ($d21 - $d9) / $d9 * 100
